# ---------------------------------------------------------------------------
# Jogos_da_Semana_FlashScore_2024-11-11.xlsx update
#
# Two logical changes are applied to Sheet1:
#
#   1. The "Odd_CS_3-3_HT" column is moved from its old position (BC) to
#      just before "Odd_CS_0-1_HT" (AW), shifting the six columns in
#      between (AW..BC) one place to the right (AX..BD). "Odd_CS_4-4_HT"
#      (BD) is unaffected - it stays the last column.
#
#   2. A brand-new match (Hobro vs Hillerod, Denmark 1st Division) is
#      inserted as the new row 2, pushing the existing match
#      (Gloria Buzau vs Petrolul, Romania Liga 1) down to row 3.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the "Odd_CS_3-3_HT" column ---------------------------
# Open a blank column at AW (xlShiftToRight = -4161), which pushes the old
# AW..BD block one column to the right (old BC's "Odd_CS_3-3_HT" data ends
# up sitting in BD).
$ws.Range("AW:AW").Insert(-4161)
# Move that relocated "Odd_CS_3-3_HT" column (now in BD) into the blank AW
# column that was just opened.
$ws.Range("BD1:BD2").Cut($ws.Range("AW1"))
# Close the now-empty BD column (xlShiftToLeft = -4161 on Delete), which
# brings the untouched "Odd_CS_4-4_HT" column back from BE into BD.
$ws.Range("BD:BD").Delete(-4161)

# --- Step 2: insert the new match as row 2 ----------------------------------
$ws.Rows.Item(2).Insert()
# Row-insert copies the header row's bold formatting onto the new blank
# row; strip it back to the default (unstyled) look used by other data rows.
$ws.Rows.Item(2).ClearFormats()

# --- Step 3: populate row 2 with the new match's data -----------------------
$ws.Cells.Item(2,1).NumberFormat = "@"   # A2 = Id
$ws.Cells.Item(2,1).Value = "G6XVkCH8"
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).NumberFormat = "@"   # B2 = Date
$ws.Cells.Item(2,2).Value = "11/11/2024"
$ws.Cells.Item(2,2).ClearFormats()
$ws.Cells.Item(2,3).NumberFormat = "@"   # C2 = Time
$ws.Cells.Item(2,3).Value = "15:00"
$ws.Cells.Item(2,3).ClearFormats()
$ws.Cells.Item(2,4).NumberFormat = "@"   # D2 = League
$ws.Cells.Item(2,4).Value = "DENMARK - 1ST DIVISION"
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).NumberFormat = "@"   # E2 = Home
$ws.Cells.Item(2,5).Value = "Hobro"
$ws.Cells.Item(2,5).ClearFormats()
$ws.Cells.Item(2,6).NumberFormat = "@"   # F2 = Away
$ws.Cells.Item(2,6).Value = "Hillerod"
$ws.Cells.Item(2,6).ClearFormats()
$ws.Cells.Item(2,7).Value = 2.5   # G2 = Odd_H_FT
$ws.Cells.Item(2,8).Value = 3.6   # H2 = Odd_D_FT
$ws.Cells.Item(2,9).Value = 2.5   # I2 = Odd_A_FT
$ws.Cells.Item(2,10).Value = 3.1   # J2 = Odd_H_HT
$ws.Cells.Item(2,11).Value = 2.3   # K2 = Odd_D_HT
$ws.Cells.Item(2,12).Value = 3.1   # L2 = Odd_A_HT
$ws.Cells.Item(2,13).Value = 1.03   # M2 = Odd_Over05_FT
$ws.Cells.Item(2,14).Value = 15   # N2 = Odd_Under05_FT
$ws.Cells.Item(2,15).Value = 1.2   # O2 = Odd_Over15_FT
$ws.Cells.Item(2,16).Value = 4.33   # P2 = Odd_Under15_FT
$ws.Cells.Item(2,17).Value = 1.65   # Q2 = Odd_Over25_FT
$ws.Cells.Item(2,18).Value = 2.2   # R2 = Odd_Under25_FT
$ws.Cells.Item(2,19).Value = 1.33   # S2 = Odd_Over05_HT
$ws.Cells.Item(2,20).Value = 3.25   # T2 = Odd_Under05_HT
$ws.Cells.Item(2,21).Value = 1.57   # U2 = Odd_BTTS_Yes
$ws.Cells.Item(2,22).Value = 2.25   # V2 = Odd_BTTS_No
$ws.Cells.Item(2,23).Value = 11   # W2 = Odd_CS_1-0
$ws.Cells.Item(2,24).Value = 15   # X2 = Odd_CS_2-0
$ws.Cells.Item(2,25).Value = 10   # Y2 = Odd_CS_2-1
$ws.Cells.Item(2,26).Value = 26   # Z2 = Odd_CS_3-0
$ws.Cells.Item(2,27).Value = 19   # AA2 = Odd_CS_3-1
$ws.Cells.Item(2,28).Value = 23   # AB2 = Odd_CS_3-2
$ws.Cells.Item(2,29).Value = 15   # AC2 = Odd_CS_0-0
$ws.Cells.Item(2,30).Value = 7.5   # AD2 = Odd_CS_1-1
$ws.Cells.Item(2,31).Value = 13   # AE2 = Odd_CS_2-2
$ws.Cells.Item(2,32).Value = 41   # AF2 = Odd_CS_3-3
$ws.Cells.Item(2,33).Value = 126   # AG2 = Odd_CS_4-4
$ws.Cells.Item(2,34).Value = 11   # AH2 = Odd_CS_0-1
$ws.Cells.Item(2,35).Value = 15   # AI2 = Odd_CS_0-2
$ws.Cells.Item(2,36).Value = 10   # AJ2 = Odd_CS_1-2
$ws.Cells.Item(2,37).Value = 26   # AK2 = Odd_CS_0-3
$ws.Cells.Item(2,38).Value = 19   # AL2 = Odd_CS_1-3
$ws.Cells.Item(2,39).Value = 23   # AM2 = Odd_CS_2-3
$ws.Cells.Item(2,40).Value = 4.75   # AN2 = Odd_CS_1-0_HT
$ws.Cells.Item(2,41).Value = 13   # AO2 = Odd_CS_2-0_HT
$ws.Cells.Item(2,42).Value = 21   # AP2 = Odd_CS_2-1_HT
$ws.Cells.Item(2,43).Value = 41   # AQ2 = Odd_CS_3-0_HT
$ws.Cells.Item(2,44).Value = 51   # AR2 = Odd_CS_3-1_HT
$ws.Cells.Item(2,45).Value = 126   # AS2 = Odd_CS_3-2_HT
$ws.Cells.Item(2,46).Value = 3.25   # AT2 = Odd_CS_0-0_HT
$ws.Cells.Item(2,47).Value = 7   # AU2 = Odd_CS_1-1_HT
$ws.Cells.Item(2,48).Value = 41   # AV2 = Odd_CS_2-2_HT
$ws.Cells.Item(2,49).Value = 351   # AW2 = Odd_CS_3-3_HT
$ws.Cells.Item(2,50).Value = 4.75   # AX2 = Odd_CS_0-1_HT
$ws.Cells.Item(2,51).Value = 13   # AY2 = Odd_CS_0-2_HT
$ws.Cells.Item(2,52).Value = 21   # AZ2 = Odd_CS_1-2_HT
$ws.Cells.Item(2,53).Value = 41   # BA2 = Odd_CS_0-3_HT
$ws.Cells.Item(2,54).Value = 51   # BB2 = Odd_CS_1-3_HT
$ws.Cells.Item(2,55).Value = 126   # BC2 = Odd_CS_2-3_HT
$ws.Cells.Item(2,56).Value = 151   # BD2 = Odd_CS_4-4_HT
